$csv = @"
,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15
0,3.473321297368443,1.027038466385278,0.04746563878543242,1.357485398655896,0.5025806726517104,0,0.01155765954546265,0.004188960627741434,0,0,0,0,0,0,0,1.409506098574781
1,3.025989452917543,0.9066745903175217,0.0426034974438565,1.182098661248517,0.4617029716188839,0,0.008202360679210235,0.002435183028116406,0,0,0,0,0,0,0,1.320823533477864
2,2.75113556242826,0.8333404639695914,0.03962876718298247,1.074759720971684,0.4372545917619561,0,0.006371091650124394,0.001625717039543506,0,0,0,0,0,0,0,1.268368109159809
3,2.63902800674191,0.8051653970680377,0.03849908900418342,1.031067517943328,0.4266486045318629,0,0.005676257758423175,0.001416686289084979,0,0,0,0,0,0,0,1.244799232327438
4,2.620371528146563,0.8024315558468516,0.03841008806588064,1.02379653887381,0.4239283188626146,0,0.005562248232508615,0.001464892024789499,0,0,0,0,0,0,0,1.237674402612328
5,2.74952406004121,0.8382659236967527,0.03988322695139601,1.074117450153068,0.4344650404707977,0,0.006356633002021828,0.001823568820628729,0,0,0,0,0,0,0,1.259198446895141
6,3.318965971132002,0.9925124729663253,0.04614815818064955,1.296854148925604,0.4848336532956665,0,0.01034268704065072,0.003760544424331158,0,0,0,0,0,0,0,1.366767437574254
7,4.435602471144477,1.290531435654884,0.05804490263918893,1.737609312871797,0.5942348880306199,0,0.02021120401831755,0.00963857218754427,0,0,0,0,0,0,0,1.6115513099856
8,5.256708081939905,1.515316710739057,0.06850537725882333,1.962773071531615,0.6648536962738518,0,0.02827314176957785,0.01570473395739747,0,0,0,0,0,0,0,1.758094049501096
9,5.621723590468036,1.619316488079619,0.08791808590716244,1.273338595727523,0.5840019648272801,0,0.0436581359854209,0.01784748950442339,0,0,0,0,0,0,0,1.458508683066697
10,5.756915616603408,1.651205380261786,0.1031250141704589,0.7762795755667042,0.5090995615705083,0,0.07943051563349712,0.01806149893065268,0,0,0,0,0,0,0,1.209200150188622
11,5.721574732849376,1.637625448490837,0.1162027654279143,0.3988614078432349,0.430380206184779,0,0.1322256592901709,0.01705535373271427,0,0,0,0,0,0,0,0.9703656183522185
12,5.619715193664376,1.60936285668663,0.124570431748225,0.2074326083771822,0.3740200338718793,0,0.1796051948735169,0.01592213451959257,0,0,0,0,0,0,0,0.8102644711294573
13,5.558434118736045,1.595017950854015,0.1260904655919148,0.1697532213726909,0.3582976480543252,0,0.1914455637791264,0.01547904891747898,0,0,0,0,0,0,0,0.7688839038803792
14,5.211112811986141,1.503985585164514,0.1185835721831836,0.1651835803463868,0.3437177782636311,0,0.1761920200204514,0.01309784018900029,0,0,0,0,0,0,0,0.7532573595373577
15,4.999469065003382,1.449566146357029,0.1083767723298337,0.2494626366359753,0.3618516828446161,0,0.1374478192042687,0.01192002461096919,0,0,0,0,0,0,0,0.8237318663030635
16,4.880056660330183,1.415658091795876,0.09497550145089662,0.4831482273805676,0.414059753865331,0,0.08531323150300807,0.01135163912682557,0,0,0,0,0,0,0,0.9919834191522341
17,4.843389989079867,1.411809549907389,0.08153301090439413,0.9182382728010765,0.4898916898786254,0,0.04284433691881162,0.01183622124126327,0,0,0,0,0,0,0,1.232335495196253
18,5.040510401192705,1.473197318833968,0.06669841540026766,1.899002168065294,0.6372030696024922,0,0.02597367173380327,0.01458617560381548,0,0,0,0,0,0,0,1.689473743140525
19,5.673470024573305,1.645288297854052,0.07226682496410319,2.232131371173438,0.7151516839723513,0,0.03394106454436541,0.02003765971659988,0,0,0,0,0,0,0,1.879691950255136
20,6.087802046524757,1.75110429291118,0.07632423049080472,2.399042766224454,0.7641198733797552,0,0.03916662059644382,0.02373310118088412,0,0,0,0,0,0,0,1.997339558723922
21,5.866700694770543,1.687950480952622,0.07382043616949829,2.309899965022154,0.7409961635137279,0,0.03634905217580897,0.02148344527393053,0,0,0,0,0,0,0,1.944593446737883
22,5.031235608746215,1.461156653618673,0.06490913351873928,1.974791371224228,0.6510398744255781,0,0.02646713233113207,0.01423350754240449,0,0,0,0,0,0,0,1.736204745891285
23,4.133146981528341,1.219123749376706,0.05531613510211741,1.617694854800973,0.5591086164494499,0,0.01727780945380286,0.00809390645974073,0,0,0,0,0,0,0,1.527106045572907
"@

$rows = $csv -split "`n"
$nrows = $rows.Count
$ncols = ($rows[0] -split ",").Count
$data = New-Object 'object[,]' $nrows,$ncols

for ($r = 0; $r -lt $nrows; $r++) {
    $fields = $rows[$r] -split ","
    for ($c = 0; $c -lt $ncols; $c++) {
        $field = $fields[$c]
        if ($field -ne "") {
            $data[$r,$c] = [double]$field
        }
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the full A1:Q25 block (header row + 24 data rows), expanding the
# table from columns A:O to A:Q (new columns P and Q).
$ws.Range("A1:Q25").Value = $data

# P1/Q1 are brand-new header cells; give them the same style (bold,
# thin border, centered) already used by the rest of row 1 (e.g. B1).
$ws.Range("B1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
